$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Películas")
$lo = $ws.ListObjects.Item(1)

# Extend the table by two rows to hold the two new movies ("Jack el caza
# gigantes" and "Mamá o papá"); the table auto re-sorts data by score, so we
# rewrite every row from the first changed row (44) through the new last row (102).
$null = $lo.ListRows.Add()
$null = $lo.ListRows.Add()

# row, movie name, Visualmente(D), Impresion personal(E), Ritmo(F), Audio(G), IMDb(H), Filmaffinity(I)
$rows = @(
    @(44, 'Jack el caza gigantes', 8, 8, 9, 8, 6.3, 5.3),
    @(45, 'Black Widow', 8, 8, 8, 8, 6.6, 5.6),
    @(46, 'Rebel Ridge', 8, 8, 8, 7, 6.8, 6),
    @(47, 'Piratas del Caribe: La venganza de Salazar', 9, 7, 7, 10, 6.5, 5.4),
    @(48, 'Flash (2023)', 7, 8, 8, 8, 6.6, 6),
    @(49, 'Proyecto Power', 7, 9, 8, 8, 6, 4.9000000000000004),
    @(50, 'Monstruos S.A.', 7, 7, 6, 7, 8.1, 7.6),
    @(51, 'Encanto', 8, 7, 7, 8, 7.2, 6.3),
    @(52, 'Bright', 7, 8, 9, 8, 6.3, 5.0999999999999996),
    @(53, 'Raya y el último dragón', 8, 7, 6, 8, 7.3, 6.6),
    @(54, 'Bichos: Una aventura en miniatura', 7, 8, 7, 6, 7.2, 6.7),
    @(55, 'Charlie y la fábrica de chocolate', 7, 7, 8, 8, 6.7, 6.6),
    @(56, 'Ocho Apellidos Vascos', 8, 8, 7, 7, 6.5, 6),
    @(57, 'Brave (Indomable)', 7, 8, 6, 7, 7.1, 6.5),
    @(58, 'A dos metros de ti', 6, 7, 8, 8, 7.2, 6.2),
    @(59, 'Spenser Confidential', 8, 8, 8, 7, 6.2, 5.2),
    @(60, 'Piratas del Caribe: En mareas misteriosas', 7, 7, 8, 9, 6.6, 5.4),
    @(61, 'Equipaje de mano', 9, 7, 8, 7, 6.5, 5.6),
    @(62, 'Mamá o papá', 8, 8, 7, 8, 6, 5.4),
    @(63, 'Liga de la Justicia (2017)', 8, 7, 8, 8, 6.1, 5.3),
    @(64, 'Lift: Un robo de primera clase', 8, 8, 8, 8, 5.5, 4.5),
    @(65, 'Ejército de ladrones', 8, 7, 8, 7, 6.4, 5.5),
    @(66, 'Lilo & Stitch (2025)', 9, 6, 6, 8, 6.9, 6.4),
    @(67, 'Alerta roja', 7, 7, 8, 8, 6.3, 5.3),
    @(68, 'El Camino: Una película de Breaking Bad', 8, 6, 7, 7, 7.3, 6.2),
    @(69, 'Shrek', 5, 6, 6, 7, 7.9, 7.8),
    @(70, 'Gozilla vs. Kong', 8, 6, 8, 7, 6.3, 5.3),
    @(71, 'Sin instrucciones', 7, 7, 7, 6, 6.3, 6.1),
    @(72, 'Chappie', 8, 6, 7, 6, 6.8, 5.8),
    @(73, 'Chip y Chop: Guardianes rescatadores', 7, 6, 6, 6, 6.9, 5.7),
    @(74, 'Paradise', 6, 6, 7, 7, 6.3, 5.6),
    @(75, 'Chang machaca', 6, 7, 6, 6, 6.4, 5.4),
    @(76, 'Morbius', 7, 6, 8, 8, 5.0999999999999996, 4.5),
    @(77, 'Cuatro Fantásticos (2015)', 8, 7, 8, 7, 4.3, 4.0999999999999996),
    @(78, 'El hoyo 2', 7, 6, 8, 8, 5, 4.2),
    @(79, 'Padre no hay más que uno 5: Nido repleto', 7, 7, 7, 7, 4.8, 4.0999999999999996),
    @(80, 'El sindicato', 7, 6, 7, 7, 5.4, 4.5999999999999996),
    @(81, 'Perdiendo el norte', 6, 6, 8, 5, 5.8, 5),
    @(82, 'Death Note (2017)', 8, 7, 6, 6, 4.5, 3.5),
    @(83, 'Red One', 7, 5, 4, 6, 6.8, 5.3),
    @(84, 'Red', 8, 3, 5, 6, 7, 6.4),
    @(85, 'K.O.', 6, 5, 6, 6, 5.8, 4.9000000000000004),
    @(86, 'Inexpertos', 3, 6, 8, 5, 5.8, 4.9000000000000004),
    @(87, 'La huella del mal', 6, 6, 6, 6, 4.8, 4.3),
    @(88, 'Estado eléctrico', 6, 4, 4, 6, 5.9, 5.0999999999999996),
    @(89, 'Barbie', 5, 3, 2, 8, 6.8, 5.8),
    @(90, 'Mi año en Oxford', 6, 3, 5, 6, 6, 4.9000000000000004),
    @(91, 'El muro negro', 5, 5, 3, 5, 5.5, 4.8),
    @(92, 'Almost cops', 5, 6, 6, 2, 4.8, 3.7),
    @(93, 'Hotel Bitcoin', 5, 5, 5, 3, 4.9000000000000004, 4),
    @(94, 'Ocho Apellidos Marroquís', 5, 4, 5, 5, 4.9000000000000004, 3.9),
    @(95, 'Los hombres lobo', 1, 4, 6, 5, 5.5, 4.4000000000000004),
    @(96, 'Lo que el viento se llevó', 7, 0, 2, 2, 8.1999999999999993, 7.9),
    @(97, 'De vuelta a la acción', 5, 3, 3, 4, 5.9, 4.9000000000000004),
    @(98, 'A descubierto', 3, 4, 4, 4, 5.4, 4.7),
    @(99, 'Alimañas', 4, 2, 6, 6, 4.9000000000000004, 4.4000000000000004),
    @(100, 'Bajo el mismo techo', 5, 4, 2, 5, 4.5, 3.9),
    @(101, 'Chicos buenos', 1, 1, 0, 1, 6.7, 5.9),
    @(102, 'Supersalidos', 0, 0, 0, 0, 7.6, 5.9)
)

foreach ($r in $rows) {
    $rowNum = $r[0]

    $bCell = $ws.Range("B$rowNum")
    $bCell.Value = $r[1]
    $bCell.HorizontalAlignment = -4131

    $dCell = $ws.Range("D$rowNum")
    $dCell.Value = $r[2]
    $dCell.HorizontalAlignment = -4108

    $eCell = $ws.Range("E$rowNum")
    $eCell.Value = $r[3]
    $eCell.HorizontalAlignment = -4108

    $fCell = $ws.Range("F$rowNum")
    $fCell.Value = $r[4]
    $fCell.HorizontalAlignment = -4108

    $gCell = $ws.Range("G$rowNum")
    $gCell.Value = $r[5]
    $gCell.HorizontalAlignment = -4108

    $hCell = $ws.Range("H$rowNum")
    $hCell.Value = $r[6]
    $hCell.NumberFormat = "0.0"
    $hCell.HorizontalAlignment = -4108

    $iCell = $ws.Range("I$rowNum")
    $iCell.Value = $r[7]
    $iCell.NumberFormat = "0.0"
    $iCell.HorizontalAlignment = -4108

    $cCell = $ws.Range("C$rowNum")
    $cCell.Formula = "=AVERAGE(D$rowNum,E$rowNum,E$rowNum,F$rowNum,G$rowNum,H$rowNum,H$rowNum,I$rowNum)"
    $cCell.NumberFormat = "0.0"
    $cCell.HorizontalAlignment = -4108
    $cCell.Font.Bold = $true
}

# Update the sheet view to reflect the edited area
$ws.Range("C102").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 86